# Scheduled market-data refresh for Shinryu_Profits workbook.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* (cols H-N) on the
# affected Leve rows across all eight crafting-job sheets. Values come from
# the latest price pull; a handful of rows gain/lose their LeveProfitHQ or
# LeveProfitNQ cell entirely depending on whether HQ/NQ craft is relevant.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3338.5806
$ws.Range("I64").Value = 3122.7693
$ws.Range("J64").Value = 3494.4443
$ws.Range("K64").Value = 3122.7693
$ws.Range("L64").Value = 3494.4443
$ws.Range("M64").Value = -2874.7693
$ws.Range("N64").Value = -3990.4443
$ws.Range("H67").Value = 3338.5806
$ws.Range("I67").Value = 3122.7693
$ws.Range("J67").Value = 3494.4443
$ws.Range("K67").Value = 3122.7693
$ws.Range("L67").Value = 3494.4443
$ws.Range("M67").Value = -2264.7693
$ws.Range("N67").Value = -5210.4443
$ws.Range("H137").Value = 7577688.5
$ws.Range("I137").Value = 11906087
$ws.Range("J137").Value = 2991.6667
$ws.Range("K137").Value = 35718261
$ws.Range("L137").Value = 8975.000100000001
$ws.Range("M137").Value = -35715711
$ws.Range("N137").Value = -14075.0001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1625.4054
$ws.Range("I45").Value = 1620.1666
$ws.Range("J45").Value = 1814
$ws.Range("K45").Value = 1620.1666
$ws.Range("L45").Value = 1814
$ws.Range("M45").Value = -1243.1666
$ws.Range("N45").Value = -2568
$ws.Range("H109").Value = 84000
$ws.Range("J109").Value = 84000
$ws.Range("L109").Value = 84000
$ws.Range("N109").Value = -86774
$ws.Range("H122").Value = 921.9
$ws.Range("I122").Value = 899
$ws.Range("J122").Value = 1013.5
$ws.Range("K122").Value = 2697
$ws.Range("L122").Value = 3040.5
$ws.Range("M122").Value = -247
$ws.Range("N122").Value = -7940.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H134").Value = 4105
$ws.Range("I134").Value = 4156.5
$ws.Range("J134").Value = 4002
$ws.Range("K134").Value = 12469.5
$ws.Range("L134").Value = 12006
$ws.Range("M134").Value = -9934.5
$ws.Range("N134").Value = -17076

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2177.6453
$ws.Range("I31").Value = 2463.2727
$ws.Range("J31").Value = 2020.55
$ws.Range("K31").Value = 2463.2727
$ws.Range("L31").Value = 2020.55
$ws.Range("M31").Value = -2168.2727
$ws.Range("N31").Value = -2610.55
$ws.Range("H34").Value = 2177.6453
$ws.Range("I34").Value = 2463.2727
$ws.Range("J34").Value = 2020.55
$ws.Range("K34").Value = 2463.2727
$ws.Range("L34").Value = 2020.55
$ws.Range("M34").Value = -2261.2727
$ws.Range("N34").Value = -2424.55
$ws.Range("H134").Value = 4088.7144
$ws.Range("I134").Value = 2221.6667
$ws.Range("K134").Value = 6665.000100000001
$ws.Range("M134").Value = -4130.000100000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1039.6487
$ws.Range("I68").Value = 934.84
$ws.Range("J68").Value = 1258
$ws.Range("K68").Value = 2804.52
$ws.Range("L68").Value = 3774
$ws.Range("M68").Value = -1993.52
$ws.Range("N68").Value = -5396
$ws.Range("H71").Value = 1039.6487
$ws.Range("I71").Value = 934.84
$ws.Range("J71").Value = 1258
$ws.Range("K71").Value = 8413.559999999999
$ws.Range("L71").Value = 11322
$ws.Range("M71").Value = -4357.559999999999
$ws.Range("N71").Value = -19434
$ws.Range("H107").Value = 35714900
$ws.Range("I107").Value = 196.52942
$ws.Range("J107").Value = 90910350
$ws.Range("K107").Value = 589.58826
$ws.Range("L107").Value = 272731050
$ws.Range("M107").Value = 1330.41174
$ws.Range("N107").Value = -272734890

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H26").Value = 14521
$ws.Range("J26").Value = 14521
$ws.Range("L26").Value = 14521
$ws.Range("N26").Value = -15081
$ws.Range("H50").Value = 14521
$ws.Range("J50").Value = 14521
$ws.Range("L50").Value = 14521
$ws.Range("N50").Value = -15517
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H80").Value = 2613.4546
$ws.Range("I80").Value = 2566
$ws.Range("J80").Value = 2715.1428
$ws.Range("K80").Value = 2566
$ws.Range("L80").Value = 2715.1428
$ws.Range("M80").Value = -1568
$ws.Range("N80").Value = -4711.1428
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H83").Value = 2613.4546
$ws.Range("I83").Value = 2566
$ws.Range("J83").Value = 2715.1428
$ws.Range("K83").Value = 12830
$ws.Range("L83").Value = 13575.714
$ws.Range("M83").Value = -7838
$ws.Range("N83").Value = -23559.714
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H86").Value = 30000
$ws.Range("J86").Value = 30000
$ws.Range("L86").Value = 30000
$ws.Range("N86").Value = -32372
$ws.Range("H89").Value = 30000
$ws.Range("J89").Value = 30000
$ws.Range("L89").Value = 90000
$ws.Range("N89").Value = -101856
$ws.Range("H102").Value = 2391.2632
$ws.Range("I102").Value = 2329.6667
$ws.Range("J102").Value = 3500
$ws.Range("K102").Value = 2329.6667
$ws.Range("L102").Value = 3500
$ws.Range("M102").Value = -707.6667000000002
$ws.Range("N102").Value = -6744
$ws.Range("H122").Value = 3204.1667
$ws.Range("I122").Value = 3380
$ws.Range("J122").Value = 3078.5715
$ws.Range("K122").Value = 10140
$ws.Range("L122").Value = 9235.7145
$ws.Range("M122").Value = -7690
$ws.Range("N122").Value = -14135.7145
$ws.Range("H126").Value = 3632.611
$ws.Range("I126").Value = 3222.077
$ws.Range("J126").Value = 4700
$ws.Range("K126").Value = 9666.231
$ws.Range("L126").Value = 14100
$ws.Range("M126").Value = -7196.231
$ws.Range("N126").Value = -19040
$ws.Range("H132").Value = 4489.3335
$ws.Range("I132").Value = 4589.9487
$ws.Range("J132").Value = 3835.3333
$ws.Range("K132").Value = 13769.8461
$ws.Range("L132").Value = 11505.9999
$ws.Range("M132").Value = -11239.8461
$ws.Range("N132").Value = -16565.9999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3672
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3672
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 11016
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -15916
$ws.Range("H132").Value = 3986.88
$ws.Range("I132").Value = 4010.0286
$ws.Range("J132").Value = 3932.8667
$ws.Range("K132").Value = 12030.0858
$ws.Range("L132").Value = 11798.6001
$ws.Range("M132").Value = -9500.085800000001
$ws.Range("N132").Value = -16858.6001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2746.7693
$ws.Range("I122").Value = 2200.4614
$ws.Range("J122").Value = 3293.077
$ws.Range("K122").Value = 6601.3842
$ws.Range("L122").Value = 9879.231
$ws.Range("M122").Value = -4151.3842
$ws.Range("N122").Value = -14779.231
$ws.Range("H132").Value = 2073.6562
$ws.Range("I132").Value = 1389.7084
$ws.Range("J132").Value = 4125.5
$ws.Range("K132").Value = 4169.1252
$ws.Range("L132").Value = 12376.5
$ws.Range("M132").Value = -1639.1252
$ws.Range("N132").Value = -17436.5
$ws.Range("H137").Value = 25128.6
$ws.Range("J137").Value = 25128.6
$ws.Range("L137").Value = 25128.6
$ws.Range("N137").Value = -35328.6
